# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows above current row 2 (which holds the first data row),
# shifting existing data down.
$insertRange = $ws.Range("A2:C8")
$insertRange.Insert()

$newTopData = @(
    @(-0.7498788833618164, 0.7729501724243164, -1.214103698730469),
    @(-0.8813939094543457, 0.8021388053894043, -1.258926272392273),
    @(-0.6222906112670898, 0.7483844757080078, -1.280380129814148),
    @(-0.7319130897521973, 0.7634215354919434, -1.176308631896973),
    @(-0.568336009979248,  0.7112784385681152, -1.411303043365478),
    @(-0.4251332283020019, 0.7228684425354004, -1.317938923835754),
    @(-0.824821949005127,  0.7479877471923828, -1.330062508583069)
)

$r = 2
foreach ($row in $newTopData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Append one new row at the bottom (row 29)
$ws.Cells.Item(29, 1).Value = -0.525787353515625
$ws.Cells.Item(29, 2).Value = 1.078789949417114
$ws.Cells.Item(29, 3).Value = -1.021092414855957
